$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44438
$ws.Range("J2").Value = 400
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 11500
$ws.Range("P2").Value = 383

$ws.Range("D3").Value = 44858
$ws.Range("H3").Value = "Española"
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 9500
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 9750
$ws.Range("N3").Value = "$/caja 30 unidades"
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 325
$ws.Range("Q3").Value = 30

$ws.Range("D4").Value = 44687
$ws.Range("H4").Value = "Española"
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 19000
$ws.Range("M4").Value = 18500
$ws.Range("N4").Value = "$/caja 30 unidades"
$ws.Range("P4").Value = 617
$ws.Range("Q4").Value = 30

$ws.Range("D5").Value = 44427
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 12500
$ws.Range("O5").Value = "Provincia de Limarí"
$ws.Range("P5").Value = 312

$ws.Range("D6").Value = 44767
$ws.Range("H6").Value = "Madrigal"
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 14500
$ws.Range("N6").Value = "$/caja 40 unidades"
$ws.Range("P6").Value = 362
$ws.Range("Q6").Value = 40

$ws.Range("D7").Value = 44484
$ws.Range("H7").Value = "Española"
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 9000
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = 9500
$ws.Range("N7").Value = "$/caja 30 unidades"
$ws.Range("O7").Value = "Provincia del Elquí"
$ws.Range("P7").Value = 317
$ws.Range("Q7").Value = 30

$ws.Range("D8").Value = 44729
$ws.Range("H8").Value = "Madrigal"
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 16000
$ws.Range("L8").Value = 17000
$ws.Range("M8").Value = 16500
$ws.Range("N8").Value = "$/caja 40 unidades"
$ws.Range("P8").Value = 412
$ws.Range("Q8").Value = 40

$ws.Range("D9").Value = 44784
$ws.Range("H9").Value = "Madrigal"
$ws.Range("J9").Value = 520
$ws.Range("K9").Value = 11500
$ws.Range("L9").Value = 12000
$ws.Range("M9").Value = 11750
$ws.Range("N9").Value = "$/caja 40 unidades"
$ws.Range("P9").Value = 294
$ws.Range("Q9").Value = 40

$ws.Range("D10").Value = 44790
$ws.Range("K10").Value = 14000
$ws.Range("L10").Value = 15000
$ws.Range("M10").Value = 14500
$ws.Range("P10").Value = 483

$ws.Range("H11").Value = "Madrigal"
$ws.Range("K11").Value = 11500
$ws.Range("L11").Value = 12000
$ws.Range("M11").Value = 11750
$ws.Range("N11").Value = "$/caja 40 unidades"
$ws.Range("O11").Value = "Provincia del Elquí"
$ws.Range("P11").Value = 294
$ws.Range("Q11").Value = 40

$ws.Range("D12").Value = 44855
$ws.Range("H12").Value = "Española"
$ws.Range("J12").Value = 540
$ws.Range("K12").Value = 9500
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = 9750
$ws.Range("N12").Value = "$/caja 30 unidades"
$ws.Range("P12").Value = 325
$ws.Range("Q12").Value = 30

$ws.Range("D13").Value = 44420
$ws.Range("J13").Value = 800
$ws.Range("K13").Value = 14000
$ws.Range("L13").Value = 15000
$ws.Range("M13").Value = 14500
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 362

$ws.Range("D14").Value = 44420
$ws.Range("H14").Value = "Madrigal"
$ws.Range("J14").Value = 700
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13500
$ws.Range("N14").Value = "$/caja 40 unidades"
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 338
$ws.Range("Q14").Value = 40

$ws.Range("D15").Value = 44701
$ws.Range("H15").Value = "Española"
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 19000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19500
$ws.Range("N15").Value = "$/caja 30 unidades"
$ws.Range("O15").Value = "Provincia del Elquí"
$ws.Range("P15").Value = 650
$ws.Range("Q15").Value = 30

$ws.Range("D16").Value = 44426
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 11500
$ws.Range("M16").Value = 11750
$ws.Range("O16").Value = "Provincia de Limarí"
$ws.Range("P16").Value = 392

$ws.Range("D17").Value = 44426
$ws.Range("H17").Value = "Madrigal"
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 12500
$ws.Range("L17").Value = 13000
$ws.Range("M17").Value = 12750
$ws.Range("N17").Value = "$/caja 40 unidades"
$ws.Range("O17").Value = "Provincia de Limarí"
$ws.Range("P17").Value = 319
$ws.Range("Q17").Value = 40

$ws.Range("D18").Value = 44498
$ws.Range("K18").Value = 8500
$ws.Range("L18").Value = 9000
$ws.Range("M18").Value = 8750
$ws.Range("O18").Value = "Provincia de Limarí"
$ws.Range("P18").Value = 292

$ws.Range("D19").Value = 44839
$ws.Range("H19").Value = "Española"
$ws.Range("J19").Value = 400
$ws.Range("K19").Value = 12000
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 12500
$ws.Range("N19").Value = "$/caja 30 unidades"
$ws.Range("P19").Value = 417
$ws.Range("Q19").Value = 30

